# Apply cryptos.xlsx price/volume update (2023-12-20 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.634.47'
$ws.Range("E2").Value = '  +3.06%  '

$ws.Range("D3").Value = '2.196.62'
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("E4").Value = '  -0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '260.27'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.08%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '82.12'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +12.44%  '

$ws.Range("E7").Value = '  +1.78%  '

$ws.Range("E8").Value = '  -0.07%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.593'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +1.96%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '43.64'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +9.58%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0919'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("E12").Value = '  +3.13%  '

$ws.Range("E13").Value = '  +2.48%  '

$ws.Range("D14").Value = '2.522.86'
$ws.Range("E14").Value = '  +0.65%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '14.27'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.75%  '

$ws.Range("D16").Value = '2.211.93'
$ws.Range("E16").Value = '  +1.55%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.780'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.46%  '

$ws.Range("D18").Value = '43.535.55'
$ws.Range("E18").Value = '  +2.99%  '

$ws.Range("E19").Value = '  +1.06%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '69.82'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.99%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '5.92'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +1.42%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '2.44'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +15.79%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '230.44'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.08%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '8.87'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -4.93%  '

$ws.Range("E25").Value = '  -0.07%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '42.28'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +14.95%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '10.68'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +2.33%  '

$ws.Range("E28").Value = '  -0.29%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.24'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +3.43%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '2.20'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +2.91%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '173.93'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +2.37%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '20.44'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.33%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0869'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +7.01%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '5.33'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.56%  '

$ws.Range("E35").Value = '  +7.10%  '

$ws.Range("E36").Value = '  +1.99%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '4.47'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +6.92%  '

$ws.Range("E38").Value = '  +5.27%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '13.12'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +12.06%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.88'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +17.98%  '

$ws.Range("E41").Value = '  +2.51%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '63.81'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +7.98%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '5.46'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +6.39%  '

$ws.Range("E44").Value = '  +2.80%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '8.27'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.54%  '

$ws.Range("E48").Value = '  +4.28%  '

$ws.Range("E49").Value = '  +2.88%  '

$ws.Range("E50").Value = '  -4.74%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.49'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +24.09%  '

# Row 45/46: Aave and Cronos swap ranking positions with updated figures
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '100.38'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0981'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.22%  '

